$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8624605536460876
$ws.Range("B1").Value = 1.716741323471069
$ws.Range("C1").Value = 2.606055736541748
$ws.Range("D1").Value = 1.531292319297791
$ws.Range("E1").Value = 0.8207937479019165
